$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new SVR parameter headers in row 1
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Add corresponding values in row 2
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Update the selected cell to match the target state
$ws.Range("K6").Select()
